# The "Recorded By" column (G) on the "Session Analysis Results" sheet lists the
# people/systems that recorded a session as a comma-separated string, e.g.
#   "System, dnasr281@gmail.com"
# Upstream flipped the display order of that list (entries are now shown in
# reverse order), EXCEPT for rows whose list already ends with the literal
# value "System" - those are left untouched (e.g. "admin@admin.com, System").
# Single-value cells (only one recorder, no comma) are also left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") { continue }

    $parts = $val -split ",\s*"
    if ($parts.Count -le 1) { continue }                 # nothing to reorder
    if ($parts[$parts.Count - 1] -eq "System") { continue }  # already ends with System -> leave as-is

    $reversed = @()
    for ($i = $parts.Count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $cell.Value = [string]::Join(", ", $reversed)
}
